# Auto-generated edit script for resum_diari_meteocat.xlsx
# Updates DATA_EXTRACCIO timestamps and several weather metric cells
# to reflect the 2026-02-24 05:49 automatic data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (outside the data range) used to push percentage-like
# text values ("NN%") into cells without Excel auto-converting them
# into numeric percentages, while preserving the destination cell's
# existing style (PasteSpecial values-only keeps the target format).
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"

# --- Plain text / numeric-with-units cell updates ---
$ws.Range("E2").Value = "2026-02-24 05:48:08"
$ws.Range("N2").Value = "0.7 °C 5:14 TU"
$ws.Range("O2").Value = "1.9 °C"
$ws.Range("E3").Value = "2026-02-24 05:48:10"
$ws.Range("N3").Value = "-0.1 °C 5:04 TU"
$ws.Range("O3").Value = "2.5 °C"
$ws.Range("E4").Value = "2026-02-24 05:48:12"
$ws.Range("J4").Value = "1022.4 hPa"
$ws.Range("O4").Value = "6.9 °C"
$ws.Range("E5").Value = "2026-02-24 05:48:15"
$ws.Range("E6").Value = "2026-02-24 05:48:17"
$ws.Range("N6").Value = "7.4 °C 5:29 TU"
$ws.Range("O6").Value = "9.4 °C"
$ws.Range("E7").Value = "2026-02-24 05:48:19"
$ws.Range("J7").Value = "1022.1 hPa"
$ws.Range("K7").Value = "-0.1 MJ/m2"
$ws.Range("N7").Value = "11.1 °C 5:02 TU"
$ws.Range("E8").Value = "2026-02-24 05:48:21"
$ws.Range("J8").Value = "1021.7 hPa"
$ws.Range("N8").Value = "13.7 °C 5:19 TU"
$ws.Range("O8").Value = "14.9 °C"
$ws.Range("E9").Value = "2026-02-24 05:48:24"
$ws.Range("E10").Value = "2026-02-24 05:48:26"
$ws.Range("E11").Value = "2026-02-24 05:48:28"
$ws.Range("N11").Value = "1.5 °C 5:26 TU"
$ws.Range("E12").Value = "2026-02-24 05:48:30"
$ws.Range("N12").Value = "3.5 °C 5:14 TU"
$ws.Range("O12").Value = "5.9 °C"
$ws.Range("E13").Value = "2026-02-24 05:48:32"
$ws.Range("O13").Value = "-1.4 °C"
$ws.Range("E14").Value = "2026-02-24 05:48:34"
$ws.Range("L14").Value = "24.1 km/h - 304º 5:27 TU"
$ws.Range("O14").Value = "9.0 °C"
$ws.Range("E15").Value = "2026-02-24 05:48:37"
$ws.Range("N15").Value = "4.0 °C 5:22 TU"
$ws.Range("O15").Value = "5.8 °C"
$ws.Range("E16").Value = "2026-02-24 05:48:38"
$ws.Range("O16").Value = "4.0 °C"
$ws.Range("E17").Value = "2026-02-24 05:48:39"
$ws.Range("L17").Value = "38.5 km/h - 232º 5:29 TU"
$ws.Range("E18").Value = "2026-02-24 05:48:40"
$ws.Range("N18").Value = "1.5 °C 5:29 TU"
$ws.Range("O18").Value = "3.0 °C"
$ws.Range("E19").Value = "2026-02-24 05:48:41"
$ws.Range("L19").Value = "7.9 km/h - 317º 5:17 TU"
$ws.Range("E20").Value = "2026-02-24 05:48:42"
$ws.Range("K20").Value = "-0.1 MJ/m2"
$ws.Range("N20").Value = "-0.7 °C 5:23 TU"
$ws.Range("O20").Value = "0.9 °C"
$ws.Range("E21").Value = "2026-02-24 05:48:43"
$ws.Range("N21").Value = "1.7 °C 5:21 TU"
$ws.Range("O21").Value = "3.5 °C"
$ws.Range("E22").Value = "2026-02-24 05:48:46"
$ws.Range("E23").Value = "2026-02-24 05:48:48"
$ws.Range("E24").Value = "2026-02-24 05:48:50"
$ws.Range("J24").Value = "1024.1 hPa"
$ws.Range("O24").Value = "3.3 °C"
$ws.Range("E25").Value = "2026-02-24 05:48:52"
$ws.Range("O25").Value = "4.8 °C"
$ws.Range("E26").Value = "2026-02-24 05:48:55"
$ws.Range("J26").Value = "1022.6 hPa"
$ws.Range("E27").Value = "2026-02-24 05:48:57"
$ws.Range("E28").Value = "2026-02-24 05:48:59"
$ws.Range("J28").Value = "1023.9 hPa"
$ws.Range("O28").Value = "3.7 °C"
$ws.Range("E29").Value = "2026-02-24 05:49:02"
$ws.Range("N29").Value = "3.5 °C 5:00 TU"
$ws.Range("O29").Value = "5.1 °C"
$ws.Range("E30").Value = "2026-02-24 05:49:04"
$ws.Range("N30").Value = "7.5 °C 5:22 TU"
$ws.Range("O30").Value = "9.1 °C"
$ws.Range("E31").Value = "2026-02-24 05:49:06"
$ws.Range("N31").Value = "13.6 °C 5:01 TU"
$ws.Range("O31").Value = "14.8 °C"
$ws.Range("E32").Value = "2026-02-24 05:49:08"
$ws.Range("E33").Value = "2026-02-24 05:49:11"
$ws.Range("N33").Value = "0.4 °C 5:02 TU"
$ws.Range("O33").Value = "2.1 °C"
$ws.Range("E34").Value = "2026-02-24 05:49:13"
$ws.Range("E35").Value = "2026-02-24 05:49:16"
$ws.Range("J35").Value = "1024.5 hPa"
$ws.Range("E36").Value = "2026-02-24 05:49:18"
$ws.Range("J36").Value = "1021.8 hPa"
$ws.Range("N36").Value = "6.8 °C 5:24 TU"
$ws.Range("O36").Value = "8.2 °C"
$ws.Range("E37").Value = "2026-02-24 05:49:20"
$ws.Range("J37").Value = "1027.5 hPa"
$ws.Range("L37").Value = "5.4 km/h - 210º 5:21 TU"
$ws.Range("N37").Value = "-0.3 °C 5:03 TU"
$ws.Range("O37").Value = "0.7 °C"
$ws.Range("E38").Value = "2026-02-24 05:49:23"
$ws.Range("N38").Value = "4.8 °C 5:29 TU"
$ws.Range("O38").Value = "6.7 °C"
$ws.Range("E39").Value = "2026-02-24 05:49:25"
$ws.Range("E40").Value = "2026-02-24 05:49:27"
$ws.Range("N40").Value = "-0.1 °C 5:28 TU"
$ws.Range("O40").Value = "1.2 °C"
$ws.Range("E41").Value = "2026-02-24 05:49:29"
$ws.Range("E42").Value = "2026-02-24 05:49:32"
$ws.Range("N42").Value = "5.2 °C 5:18 TU"
$ws.Range("O42").Value = "6.8 °C"
$ws.Range("E43").Value = "2026-02-24 05:49:34"
$ws.Range("O43").Value = "4.1 °C"
$ws.Range("E44").Value = "2026-02-24 05:49:36"
$ws.Range("L44").Value = "9.4 km/h - 9º 5:29 TU"
$ws.Range("N44").Value = "-3.2 °C 5:23 TU"
$ws.Range("O44").Value = "-0.2 °C"
$ws.Range("E45").Value = "2026-02-24 05:49:38"
$ws.Range("J45").Value = "1025.5 hPa"
$ws.Range("O45").Value = "4.5 °C"
$ws.Range("E46").Value = "2026-02-24 05:49:41"
$ws.Range("J46").Value = "1023.9 hPa"
$ws.Range("N46").Value = "0.9 °C 5:01 TU"
$ws.Range("O46").Value = "2.3 °C"

# --- Percentage text cell updates (kept as text, e.g. "88%") ---
$helper.Value = "88%"
$helper.Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4163) | Out-Null
$helper.Value = "35%"
$helper.Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4163) | Out-Null
$helper.Value = "81%"
$helper.Copy() | Out-Null
$ws.Range("H7").PasteSpecial(-4163) | Out-Null
$helper.Value = "94%"
$helper.Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4163) | Out-Null
$helper.Value = "90%"
$helper.Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4163) | Out-Null
$helper.Value = "19%"
$helper.Copy() | Out-Null
$ws.Range("H16").PasteSpecial(-4163) | Out-Null
$helper.Value = "31%"
$helper.Copy() | Out-Null
$ws.Range("H17").PasteSpecial(-4163) | Out-Null
$helper.Value = "17%"
$helper.Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4163) | Out-Null
$helper.Value = "27%"
$helper.Copy() | Out-Null
$ws.Range("H23").PasteSpecial(-4163) | Out-Null
$helper.Value = "54%"
$helper.Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4163) | Out-Null
$helper.Value = "38%"
$helper.Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4163) | Out-Null
$helper.Value = "83%"
$helper.Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4163) | Out-Null
$helper.Value = "52%"
$helper.Copy() | Out-Null
$ws.Range("H31").PasteSpecial(-4163) | Out-Null
$helper.Value = "55%"
$helper.Copy() | Out-Null
$ws.Range("H34").PasteSpecial(-4163) | Out-Null
$helper.Value = "97%"
$helper.Copy() | Out-Null
$ws.Range("H36").PasteSpecial(-4163) | Out-Null
$helper.Value = "95%"
$helper.Copy() | Out-Null
$ws.Range("H40").PasteSpecial(-4163) | Out-Null
$helper.Value = "78%"
$helper.Copy() | Out-Null
$ws.Range("H41").PasteSpecial(-4163) | Out-Null
$helper.Value = "60%"
$helper.Copy() | Out-Null
$ws.Range("H45").PasteSpecial(-4163) | Out-Null
$helper.Value = "99%"
$helper.Copy() | Out-Null
$ws.Range("H46").PasteSpecial(-4163) | Out-Null

# Clean up the helper cell so it does not leave stray data behind.
$helper.Clear() | Out-Null
$excel.CutCopyMode = 0

